$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2646.2
$ws.Range("J17").Value = 2646.2
$ws.Range("L17").Value = 7938.599999999999
$ws.Range("N17").Value = -8274.599999999999
$ws.Range("H31").Value = 2975
$ws.Range("I31").Value = 2975
$ws.Range("K31").Value = 8925
$ws.Range("M31").Value = -8695
$ws.Range("H32").Value = 927.1429000000001
$ws.Range("I32").Value = 747.5
$ws.Range("J32").Value = 999
$ws.Range("K32").Value = 747.5
$ws.Range("L32").Value = 999
$ws.Range("M32").Value = -421.5
$ws.Range("N32").Value = -1651
$ws.Range("H113").Value = 17674.5
$ws.Range("I113").Value = 13449.5
$ws.Range("J113").Value = 21899.5
$ws.Range("K113").Value = 13449.5
$ws.Range("L113").Value = 21899.5
$ws.Range("M113").Value = -10195.5
$ws.Range("N113").Value = -28407.5
$ws.Range("H132").Value = 2261.96
$ws.Range("I132").Value = 2147.875
$ws.Range("K132").Value = 6443.625
$ws.Range("M132").Value = -3913.625
$ws.Range("H135").Value = 891
$ws.Range("I135").Value = 891
$ws.Range("K135").Value = 8019
$ws.Range("M135").Value = -5484
$ws.Range("H138").Value = 2510.742
$ws.Range("I138").Value = 1371.3334
$ws.Range("K138").Value = 4114.0002
$ws.Range("M138").Value = 1025.9998

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H94").Value = 64582.25
$ws.Range("J94").Value = 64582.25
$ws.Range("L94").Value = 64582.25
$ws.Range("N94").Value = -66384.25
$ws.Range("H110").Value = 1350.7778
$ws.Range("I110").Value = 1350.7778
$ws.Range("K110").Value = 1350.7778
$ws.Range("M110").Value = 694.2221999999999
$ws.Range("H132").Value = 815.5
$ws.Range("J132").Value = 588
$ws.Range("L132").Value = 1764
$ws.Range("N132").Value = -6824

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1384.1666
$ws.Range("J5").Value = 3005
$ws.Range("L5").Value = 3005
$ws.Range("N5").Value = -3231
$ws.Range("H7").Value = 857263.5600000001
$ws.Range("I7").Value = 1200034
$ws.Range("J7").Value = 337.5
$ws.Range("K7").Value = 1200034
$ws.Range("L7").Value = 337.5
$ws.Range("M7").Value = -1199921
$ws.Range("N7").Value = -563.5
$ws.Range("H76").Value = 14846.286
$ws.Range("J76").Value = 14846.286
$ws.Range("L76").Value = 14846.286
$ws.Range("N76").Value = -15476.286
$ws.Range("H79").Value = 14846.286
$ws.Range("J79").Value = 14846.286
$ws.Range("L79").Value = 14846.286
$ws.Range("N79").Value = -17030.286
$ws.Range("H105").Value = 3403.4614
$ws.Range("I105").Value = 3260.7778
$ws.Range("J105").Value = 3724.5
$ws.Range("K105").Value = 3260.7778
$ws.Range("L105").Value = 3724.5
$ws.Range("M105").Value = -1513.7778
$ws.Range("N105").Value = -7218.5
$ws.Range("H107").Value = 1125.5834
$ws.Range("I107").Value = 950.8
$ws.Range("J107").Value = 1999.5
$ws.Range("K107").Value = 950.8
$ws.Range("L107").Value = 1999.5
$ws.Range("M107").Value = 969.2
$ws.Range("N107").Value = -5839.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2991.875
$ws.Range("I31").Value = 2131.7
$ws.Range("J31").Value = 4425.5
$ws.Range("K31").Value = 2131.7
$ws.Range("L31").Value = 4425.5
$ws.Range("M31").Value = -1836.7
$ws.Range("N31").Value = -5015.5
$ws.Range("H34").Value = 2991.875
$ws.Range("I34").Value = 2131.7
$ws.Range("J34").Value = 4425.5
$ws.Range("K34").Value = 2131.7
$ws.Range("L34").Value = 4425.5
$ws.Range("M34").Value = -1929.7
$ws.Range("N34").Value = -4829.5
$ws.Range("H58").Value = 5307.4116
$ws.Range("I58").Value = 4331.2144
$ws.Range("K58").Value = 4331.2144
$ws.Range("M58").Value = -4128.2144
$ws.Range("H134").Value = 2576.75
$ws.Range("I134").Value = 2109.75
$ws.Range("K134").Value = 6329.25
$ws.Range("M134").Value = -3794.25
$ws.Range("H136").Value = 5307.4116
$ws.Range("I136").Value = 4331.2144
$ws.Range("K136").Value = 12993.6432
$ws.Range("M136").Value = -10443.6432

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 312.69232
$ws.Range("I7").Value = 226.66667
$ws.Range("J7").Value = 386.42856
$ws.Range("K7").Value = 680.00001
$ws.Range("L7").Value = 1159.28568
$ws.Range("M7").Value = -568.00001
$ws.Range("N7").Value = -1383.28568
$ws.Range("H11").Value = 334451.16
$ws.Range("I11").Value = 667900.3
$ws.Range("K11").Value = 2003700.9
$ws.Range("M11").Value = -2003560.9
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = $null
$ws.Range("M33").Value = $null
$ws.Range("N33").Value = 0
$ws.Range("H34").Value = 527.2
$ws.Range("I34").Value = 610
$ws.Range("J34").Value = 403
$ws.Range("K34").Value = 1830
$ws.Range("L34").Value = 1209
$ws.Range("M34").Value = -1746
$ws.Range("N34").Value = -1377
$ws.Range("H122").Value = 950
$ws.Range("I122").Value = 950
$ws.Range("K122").Value = 8550
$ws.Range("M122").Value = -6100
$ws.Range("H139").Value = 1224.1875
$ws.Range("J139").Value = 1265.8
$ws.Range("L139").Value = 3797.4
$ws.Range("N139").Value = -14077.4

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1098.4
$ws.Range("I5").Value = 1098.4
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1098.4
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = $null
$ws.Range("N5").Value = -986.4000000000001
$ws.Range("H132").Value = 3554
$ws.Range("I132").Value = 3831
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 11493
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -8963
$ws.Range("N132").Value = -14060

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = $null
$ws.Range("M22").Value = $null
$ws.Range("N22").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = $null
$ws.Range("M27").Value = $null
$ws.Range("N27").Value = 0
$ws.Range("H46").Value = 4000
$ws.Range("I46").Value = 2500
$ws.Range("J46").Value = 4600
$ws.Range("K46").Value = 2500
$ws.Range("L46").Value = 4600
$ws.Range("M46").Value = -2312
$ws.Range("N46").Value = -4976

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 7999.4
$ws.Range("H4").Value = 4400.8
$ws.Range("I4").Value = 502
$ws.Range("J4").Value = 7000
$ws.Range("K4").Value = 502
$ws.Range("L4").Value = 7000
$ws.Range("M4").Value = -389
$ws.Range("N4").Value = -7226
$ws.Range("H100").Value = 947.63635
$ws.Range("I100").Value = 1010.55554
$ws.Range("K100").Value = 2021.11108
$ws.Range("M100").Value = -1480.11108
$ws.Range("H109").Value = 91584.664
$ws.Range("J109").Value = 91584.664
$ws.Range("L109").Value = 91584.664
$ws.Range("N109").Value = -94358.664
$ws.Range("H122").Value = 3004
$ws.Range("I122").Value = 3004
$ws.Range("K122").Value = 9012
$ws.Range("M122").Value = -6562
$ws.Range("H132").Value = 1382.0769
$ws.Range("I132").Value = 1379.2727
$ws.Range("K132").Value = 4137.8181
$ws.Range("M132").Value = -1607.8181
